$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 11
$ws.Range("D2").Value = 0.606
$ws.Range("E2").Value = 0.805
$ws.Range("F2").Value = 0.798
$ws.Range("G2").Value = 0.39
$ws.Range("H2").Value = 1.528

# Row 3
$ws.Range("D3").Value = 0.606
$ws.Range("E3").Value = 1.176
$ws.Range("F3").Value = 1.151
$ws.Range("G3").Value = 0.195
$ws.Range("H3").Value = 1.848

# Row 4
$ws.Range("C4").Value = 8
$ws.Range("D4").Value = 0.613
$ws.Range("E4").Value = 0.984
$ws.Range("F4").Value = 1.037
$ws.Range("G4").Value = 0.262
$ws.Range("H4").Value = 1.609

# Row 5
$ws.Range("D5").Value = 0.618
$ws.Range("E5").Value = 1.141
$ws.Range("F5").Value = 1.328
$ws.Range("G5").Value = 0.257
$ws.Range("H5").Value = 1.949
